# Auto-update draw results: append the 2025-10-18 Pick 3 draw as a new
# row (row 32) at the bottom of the results table on the active sheet.
#
# All columns in this table are stored as literal text (dates, the
# zero-padded "phase" code, and the dash-separated result are NOT real
# numbers/dates), so the destination cells are pre-formatted as Text
# ("@") before the values are written. This stops Excel's automatic
# type-inference from turning "2025-10-18" into a date serial or
# "251018" into a number, which would silently corrupt the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 32

$ws.Range("A${newRow}:E${newRow}").NumberFormat = "@"

$ws.Range("A${newRow}").Value = "2025-10-18"
$ws.Range("B${newRow}").Value = "Pick 3"
$ws.Range("C${newRow}").Value = "251018"
$ws.Range("D${newRow}").Value = "9-1-7"
$ws.Range("E${newRow}").Value = "2025-10-18T21:35:22.965+04:00"
